# Weekly update: insert one new daily/weekly price record for
# "Hortaliza, Vega Modelo de Temuco - Acelga" at row 130, pushing the
# existing rows 130-248 down to 131-249 (matches commit message:
# "Fruta / hortaliza, semanal").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 130, shifting 130..248 -> 131..249.
$ws.Rows.Item(130).Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Range("A130").Value = 10
$ws.Range("B130").Value = "Vega Modelo de Temuco"
$ws.Range("C130").Value = "La Araucanía"
$ws.Range("D130").Value = 44566
$ws.Range("E130").Value = 9
$ws.Range("F130").Value = 100112009
$ws.Range("G130").Value = "Acelga"
$ws.Range("H130").Value = "Sin especificar"
$ws.Range("I130").Value = "Primera"
$ws.Range("J130").Value = 20
$ws.Range("K130").Value = 8000
$ws.Range("L130").Value = 8000
$ws.Range("M130").Value = 8000
$ws.Range("N130").Value = "$/docena de atados (12 kilos)"
$ws.Range("O130").Value = "Provincia de Cautín"
$ws.Range("P130").Value = 667
$ws.Range("Q130").Value = 12
$ws.Range("R130").Value = "Hortaliza"
